# Re-run of the analysis scripts: the ml_results export now also reports a
# confidence interval for every estimate, so two new columns - "ci.lower"
# and "ci.upper" - are appended after the existing "p.value" column (G:H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers
$ws.Range("G1").Value2 = "ci.lower"
$ws.Range("H1").Value2 = "ci.upper"

# Rows that carry a confidence interval (rows 12-14 have none, matching the
# source data - those rows never had SE/df/t.ratio/p.value either, so G/H
# stay blank there too).
$rows = @(2,3,4,5,6,7,8,9,10,11,15,16,17,18,19,20,21,22,23,24,25,26)

$lower = @(
    0.243473756552374,
    0.180382487516955,
    0.140118032582772,
    0.42319506019761,
    0.383690008301359,
    0.455842220732477,
    0.354090444235573,
    0.16087773066909,
    0.039708414648435,
    0.027252348426535,
    -0.262525999583034,
    -0.663426013784709,
    -0.785739657806024,
    -0.578534979797882,
    -1.46200816716034,
    -1.58422215263512,
    0.027252348426535,
    0.321755461338179,
    0.0688690526332789,
    0.813104010528805,
    -0.0205525545345299,
    0.810069976030943
)

$upper = @(
    0.63117579494439,
    0.346903747821977,
    0.289267489898941,
    0.813870317807316,
    0.792111076317558,
    0.876655915795182,
    0.731004083580168,
    0.317458148241233,
    0.0835468938214413,
    0.0706483644306839,
    -0.11790916647007,
    -0.297966709675072,
    -0.359429280119888,
    -0.280236065165543,
    -0.708180888471145,
    -0.767380016602718,
    0.0706483644306839,
    0.634916296482465,
    0.17853455607864,
    1.60448865381459,
    0.0818168479205395,
    1.6027964865929
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 7).Value2 = $lower[$i]
    $ws.Cells.Item($r, 8).Value2 = $upper[$i]
}
